# Apply quarterly update: insert two new quarter columns (D, E) before the
# existing data (old D:K shifts right to F:M), then populate the two new
# columns with the new quarter figures (periods ending 2019-01-31 and
# 2018-10-31) across all three statements (Income Statement, Balance Sheet,
# Cash Flow Statement).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPAR")

# Insert two blank columns at D:E; existing D:K data shifts right to F:M.
$ws.Columns("D:E").Insert()

# Carry the number formatting across from the (now-shifted) old column D
# so the new D:E columns look like the rest of the table (date row uses the
# date format, data rows use the numeric format) instead of Excel's default.
# Done per contiguous data block so the blank separator rows (36/37, 78/79)
# between the three statements aren't touched.
$ws.Range("F7:M35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:M77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:M102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns with their figures.
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 233000
$ws.Range("E8").Value2 = 226200
$ws.Range("D9").Value2 = 209200
$ws.Range("E9").Value2 = 200000
$ws.Range("D10").Value2 = 23900
$ws.Range("E10").Value2 = 26200
$ws.Range("D12").Value2 = 2000
$ws.Range("E12").Value2 = 2100
$ws.Range("D13:E13").Value2 = 0
$ws.Range("D14").Value2 = 400
$ws.Range("E14").Value2 = 500
$ws.Range("D15:E15").Value2 = 0
$ws.Range("D17").Value2 = 230800
$ws.Range("E17").Value2 = 219800
$ws.Range("D18").Value2 = 2200
$ws.Range("E18").Value2 = 6400
$ws.Range("D20").Value2 = -400
$ws.Range("E20").Value2 = 100
$ws.Range("D21").Value2 = 4600
$ws.Range("E21").Value2 = 9100
$ws.Range("D22").Value2 = 300
$ws.Range("E22").Value2 = 200
$ws.Range("D23").Value2 = 1600
$ws.Range("E23").Value2 = 6300
$ws.Range("D24").Value2 = -300
$ws.Range("E24").Value2 = 1000
$ws.Range("D25:E25").Value2 = 0
$ws.Range("D26").Value2 = 1800
$ws.Range("E26").Value2 = 5200
$ws.Range("D27").Value2 = 1800
$ws.Range("E27").Value2 = 5200
$ws.Range("D28:E28").Value2 = 0
$ws.Range("D29:E29").Value2 = "NA"
$ws.Range("D30:E30").Value2 = 0
$ws.Range("D31:E31").Value2 = 0
$ws.Range("D32").Value2 = 400
$ws.Range("E32").Value2 = -100
$ws.Range("D33").Value2 = 1800
$ws.Range("E33").Value2 = 5200
$ws.Range("D34:E34").Value2 = 0
$ws.Range("D35").Value2 = 1800
$ws.Range("E35").Value2 = 5200
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 27400
$ws.Range("E41").Value2 = 15700
$ws.Range("D42:E42").Value2 = 0
$ws.Range("D43").Value2 = 142800
$ws.Range("E43").Value2 = 153500
$ws.Range("D44").Value2 = 70000
$ws.Range("E44").Value2 = 75800
$ws.Range("D45").Value2 = 5100
$ws.Range("E45").Value2 = 4800
$ws.Range("D46").Value2 = 245300
$ws.Range("E46").Value2 = 249800
$ws.Range("D47:E47").Value2 = 0
$ws.Range("D48").Value2 = 56600
$ws.Range("E48").Value2 = 55500
$ws.Range("D49").Value2 = 42400
$ws.Range("E49").Value2 = 36200
$ws.Range("D50:E50").Value2 = 0
$ws.Range("D51:E51").Value2 = 0
$ws.Range("D52").Value2 = 9500
$ws.Range("E52").Value2 = 8300
$ws.Range("D53:E53").Value2 = 0
$ws.Range("D54").Value2 = 353800
$ws.Range("E54").Value2 = 349900
$ws.Range("D57").Value2 = 76400
$ws.Range("E57").Value2 = 92600
$ws.Range("D58:E58").Value2 = 100
$ws.Range("D59").Value2 = 61600
$ws.Range("E59").Value2 = 50500
$ws.Range("D60").Value2 = 138100
$ws.Range("E60").Value2 = 143200
$ws.Range("D61").Value2 = 25500
$ws.Range("E61").Value2 = 18600
$ws.Range("D62").Value2 = 4100
$ws.Range("E62").Value2 = 4500
$ws.Range("D63:E63").Value2 = 0
$ws.Range("D64:E64").Value2 = 0
$ws.Range("D65:E65").Value2 = 0
$ws.Range("D66").Value2 = 167000
$ws.Range("E66").Value2 = 165500
$ws.Range("D68:E68").Value2 = 0
$ws.Range("D69:E69").Value2 = 0
$ws.Range("D70:E70").Value2 = 0
$ws.Range("D71:E71").Value2 = 0
$ws.Range("D72").Value2 = 103600
$ws.Range("E72").Value2 = 103900
$ws.Range("D73:E73").Value2 = 0
$ws.Range("D74:E74").Value2 = 0
$ws.Range("D75:E75").Value2 = 0
$ws.Range("D76").Value2 = 186700
$ws.Range("E76").Value2 = 184400
$ws.Range("D77:E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 1800
$ws.Range("E81").Value2 = 5200
$ws.Range("D83").Value2 = 2700
$ws.Range("E83").Value2 = 2600
$ws.Range("D84:E84").Value2 = 0
$ws.Range("D85:E85").Value2 = 0
$ws.Range("D86:E86").Value2 = 0
$ws.Range("D87:E87").Value2 = 0
$ws.Range("D88:E88").Value2 = 0
$ws.Range("D89").Value2 = 14000
$ws.Range("E89").Value2 = -3300
$ws.Range("D91").Value2 = -1600
$ws.Range("E91").Value2 = -3300
$ws.Range("D92:E92").Value2 = 0
$ws.Range("D93:E93").Value2 = 0
$ws.Range("D94").Value2 = -6800
$ws.Range("E94").Value2 = -3300
$ws.Range("D96").Value2 = -1800
$ws.Range("E96").Value2 = 0
$ws.Range("D97:E97").Value2 = 0
$ws.Range("D98:E98").Value2 = 0
$ws.Range("D99:E99").Value2 = 0
$ws.Range("D100").Value2 = 4600
$ws.Range("E100").Value2 = 700
$ws.Range("D101:E101").Value2 = 0
$ws.Range("D102").Value2 = 11800
$ws.Range("E102").Value2 = -6000
